# Updated cryptos list on Mon Nov  4 04:51:27 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original text formatting so values
# like "10.80" or "339.50" are not silently converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.192.06"
$ws.Range("E2").Value = "  +0.92%  "
$ws.Range("D3").Value = "2.475.19"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "561.42"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("D6").Value = "163.16"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").Value = "2.475.21"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  +0.80%  "
$ws.Range("D12").Value = "4.89"
$ws.Range("E12").Value = "  +1.63%  "
$ws.Range("E13").Value = "  -3.08%  "
$ws.Range("D14").Value = "69.064.55"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "23.76"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "2.476.58"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("D19").Value = "10.80"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "339.50"
$ws.Range("E20").Value = "  -3.70%  "
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "1.89"
$ws.Range("E23").Value = "  +1.06%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "67.30"
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("E26").Value = "  -1.71%  "
$ws.Range("D27").Value = "2.604.00"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").Value = "0.0₃0825"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("E31").Value = "  -1.64%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "432.81"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("D36").Value = "157.74"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").Value = "17.84"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "0.302"
$ws.Range("E41").Value = "  -1.59%  "
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("E43").Value = "  -3.65%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "133.17"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "0.486"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "0.564"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("E51").Value = "  +0.16%  "

Write-Output "Applied cryptos update."
